# Applies the crypto price/volume update described in the commit diff.
# For each changed row, the affected cells in columns B-E are updated.
# NumberFormat is forced to Text ("@") before writing the Price/Volume
# values so strings such as "0.999" or "88.452.82" are preserved as text
# (matching the original inline-string cells) instead of becoming numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "88.452.82"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "3.133.61"
$ws.Range("E3").Value = "  -4.11%  "

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "208.07"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "610.13"
$ws.Range("E6").Value = "  -2.37%  "

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.376"
$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.676"
$ws.Range("E8").Value = "  -5.25%  "

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "3.132.89"
$ws.Range("E10").Value = "  -4.01%  "

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.563"
$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.83%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.16%  "

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "88.415.76"
$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "3.708.40"
$ws.Range("E15").Value = "  -4.06%  "

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "5.16"
$ws.Range("E16").Value = "  -5.36%  "

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("D17").Value = "31.99"
$ws.Range("E17").Value = "  -5.70%  "

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.165.70"
$ws.Range("E18").Value = "  -3.89%  "

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "3.19"
$ws.Range("E19").Value = "  +1.45%  "

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "13.18"
$ws.Range("E20").Value = "  -5.97%  "

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "429.98"
$ws.Range("E21").Value = "  -1.20%  "

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "8.40"
$ws.Range("E22").Value = "  -5.04%  "

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000181"
$ws.Range("E23").Value = "  +34.61%  "

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "5.00"
$ws.Range("E24").Value = "  -5.91%  "

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "5.01"
$ws.Range("E25").Value = "  -3.84%  "

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("E26").Value = "  -5.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.353.70"

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "73.90"
$ws.Range("E28").Value = "  -3.72%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30:E30").NumberFormat = "@"
$ws.Range("D30").Value = "0.163"
$ws.Range("E30").Value = "  -9.08%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32:E32").NumberFormat = "@"
$ws.Range("D32").Value = "3.95"
$ws.Range("E32").Value = "  +28.26%  "

$ws.Range("D33:E33").NumberFormat = "@"
$ws.Range("D33").Value = "8.28"
$ws.Range("E33").Value = "  -4.62%  "

$ws.Range("D34:E34").NumberFormat = "@"
$ws.Range("D34").Value = "522.12"
$ws.Range("E34").Value = "  -6.81%  "

$ws.Range("D35:E35").NumberFormat = "@"
$ws.Range("D35").Value = "6.85"
$ws.Range("E35").Value = "  -3.35%  "

$ws.Range("D36:E36").NumberFormat = "@"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -6.21%  "

$ws.Range("D37:E37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  -8.91%  "

$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "22.21"
$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "21.54"
$ws.Range("E39").Value = "  -4.69%  "

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.995"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.125"
$ws.Range("E41").Value = "  -10.64%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("B43:E43").NumberFormat = "@"
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").Value = "  -7.04%  "

$ws.Range("B44:E44").NumberFormat = "@"
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "0.365"
$ws.Range("E44").Value = "  -8.82%  "

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "150.78"
$ws.Range("E45").Value = "  -3.04%  "

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "43.72"
$ws.Range("E46").Value = "  -2.22%  "

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = "168.87"
$ws.Range("E47").Value = "  -6.08%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.87%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.50%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.09%  "

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "3.97"
$ws.Range("E51").Value = "  -5.45%  "
